$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAP-Con")
$ws.Range("B20:B25").Insert(-4161)
